$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36 corresponds to Idaho, which previously failed with a timeout error.
# This run succeeded and produced the following values.

$ws.Range("B36").Value = 44022
$ws.Range("B36").NumberFormat = "YYYY-MM-DD"

$ws.Range("C36").Value = 9928
$ws.Range("D36").Value = 101
$ws.Range("E36").Value = 145
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 1.46
$ws.Range("H36").Value = 0.99

$ws.Range("I36").Value = $false
$ws.Range("J36").Value = $true

$ws.Range("O36").Value = "Success!"
